$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 63 (existing rows 63:73 shift down to 64:74)
$ws.Rows.Item(63).Insert()

# Fill the columns that stay constant throughout the table
$ws.Range("A63").Value = 11
$ws.Range("B63").Value = "Vega Monumental Concepción"
$ws.Range("C63").Value = "Bíobío"
$ws.Range("E63").Value = 8
$ws.Range("F63").Value = 100112024
$ws.Range("G63").Value = "Choclo"
$ws.Range("R63").Value = "Hortaliza"

# Fill the new record's specific data
$ws.Range("D63").Value = 44588
$ws.Range("H63").Value = "Choclero"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 27000
$ws.Range("K63").Value = 150
$ws.Range("L63").Value = 200
$ws.Range("M63").Value = 178
$ws.Range("N63").Value = "$/unidad"
$ws.Range("O63").Value = "Región de O'Higgins"
$ws.Range("P63").Value = 178
$ws.Range("Q63").Value = 1
